# Update the SQL queries in cells B2, C2, B3, B4, B5, B6, B7 on Sheet1 to
# use the renamed join columns (study_id / participant_id instead of id),
# and widen column C to match the longer query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToUpdate = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellRef in $cellsToUpdate) {
    $cell = $ws.Range($cellRef)
    $text = $cell.Value2

    $text = $text.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $text = $text.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $cell.Value2 = $text
}

# Widen column C to fit the now-longer query text (and drop the bestFit
# auto-sizing that Excel previously computed). 68 is the closest achievable
# ColumnWidth value to the target stored width of 68.83203125 given Excel's
# internal pixel-based column width quantization.
$ws.Columns.Item(3).ColumnWidth = 68
